$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1847.5
$ws.Range("I40").Value = 1591.6666
$ws.Range("J40").Value = 2231.25
$ws.Range("K40").Value = 1591.6666
$ws.Range("L40").Value = 2231.25
$ws.Range("M40").Value = -1416.6666
$ws.Range("N40").Value = -2581.25
$ws.Range("H43").Value = 843.7857
$ws.Range("I43").Value = 674.625
$ws.Range("J43").Value = 1069.3334
$ws.Range("K43").Value = 674.625
$ws.Range("L43").Value = 1069.3334
$ws.Range("M43").Value = -605.625
$ws.Range("N43").Value = -1207.3334
$ws.Range("H98").Value = 1518.1177
$ws.Range("I98").Value = 991.1579
$ws.Range("J98").Value = 2185.6
$ws.Range("K98").Value = 991.1579
$ws.Range("L98").Value = 2185.6
$ws.Range("M98").Value = 506.8421
$ws.Range("N98").Value = -5181.6
$ws.Range("H116").Value = 5137.25
$ws.Range("I116").Value = 1954.091
$ws.Range("J116").Value = 12140.2
$ws.Range("K116").Value = 1954.091
$ws.Range("L116").Value = 12140.2
$ws.Range("M116").Value = 1487.909
$ws.Range("N116").Value = -19024.2
$ws.Range("H122").Value = 1518.1177
$ws.Range("I122").Value = 991.1579
$ws.Range("J122").Value = 2185.6
$ws.Range("K122").Value = 2973.4737
$ws.Range("L122").Value = 6556.799999999999
$ws.Range("M122").Value = -523.4737
$ws.Range("N122").Value = -11456.8
$ws.Range("H138").Value = 2704902
$ws.Range("I138").Value = 1231.6459
$ws.Range("J138").Value = 7696293
$ws.Range("K138").Value = 3694.9377
$ws.Range("L138").Value = 23088879
$ws.Range("M138").Value = 1445.0623
$ws.Range("N138").Value = -23099159

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 27969.316
$ws.Range("I74").Value = 41065.4
$ws.Range("K74").Value = 41065.4
$ws.Range("M74").Value = -40191.4
$ws.Range("H77").Value = 27969.316
$ws.Range("I77").Value = 41065.4
$ws.Range("K77").Value = 205327
$ws.Range("M77").Value = -200959
$ws.Range("H110").Value = 1512
$ws.Range("I110").Value = 1458.8948
$ws.Range("J110").Value = 1656.1428
$ws.Range("K110").Value = 1458.8948
$ws.Range("L110").Value = 1656.1428
$ws.Range("M110").Value = 586.1052
$ws.Range("N110").Value = -5746.1428
$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29169.5
$ws.Range("J55").Value = 29169.5
$ws.Range("L55").Value = 29169.5
$ws.Range("N55").Value = -29715.5
$ws.Range("H134").Value = 803918
$ws.Range("I134").Value = 1179094
$ws.Range("J134").Value = 6668.9375
$ws.Range("K134").Value = 3537282
$ws.Range("L134").Value = 20006.8125
$ws.Range("M134").Value = -3534747
$ws.Range("N134").Value = -25076.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3529.45
$ws.Range("I58").Value = 3694
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 3694
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = -3491
$ws.Range("N58").Value = -1906
$ws.Range("H62").Value = 2365.25
$ws.Range("I62").Value = 2279.5833
$ws.Range("J62").Value = 2493.75
$ws.Range("K62").Value = 2279.5833
$ws.Range("L62").Value = 2493.75
$ws.Range("M62").Value = -1655.5833
$ws.Range("N62").Value = -3741.75
$ws.Range("H65").Value = 2365.25
$ws.Range("I65").Value = 2279.5833
$ws.Range("J65").Value = 2493.75
$ws.Range("K65").Value = 11397.9165
$ws.Range("L65").Value = 12468.75
$ws.Range("M65").Value = -8277.916499999999
$ws.Range("N65").Value = -18708.75
$ws.Range("H107").Value = 1620.0714
$ws.Range("I107").Value = 567.55
$ws.Range("J107").Value = 4251.375
$ws.Range("K107").Value = 567.55
$ws.Range("L107").Value = 4251.375
$ws.Range("M107").Value = 1352.45
$ws.Range("N107").Value = -8091.375
$ws.Range("H136").Value = 3529.45
$ws.Range("I136").Value = 3694
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 11082
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -8532
$ws.Range("N136").Value = -9600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 194.57143
$ws.Range("I12").Value = 263.1
$ws.Range("J12").Value = 156.5
$ws.Range("K12").Value = 789.3000000000001
$ws.Range("L12").Value = 469.5
$ws.Range("M12").Value = -616.3000000000001
$ws.Range("N12").Value = -815.5
$ws.Range("H22").Value = 3650
$ws.Range("J22").Value = 4220
$ws.Range("L22").Value = 12660
$ws.Range("N22").Value = -12998
$ws.Range("H27").Value = 3650
$ws.Range("J27").Value = 4220
$ws.Range("L27").Value = 12660
$ws.Range("N27").Value = -12864
$ws.Range("H113").Value = 1212584.6
$ws.Range("I113").Value = 1515650.4
$ws.Range("J113").Value = 322
$ws.Range("K113").Value = 4546951.199999999
$ws.Range("L113").Value = 966
$ws.Range("M113").Value = -4544781.199999999
$ws.Range("N113").Value = -5306
$ws.Range("H131").Value = 911.9
$ws.Range("I131").Value = 782.5
$ws.Range("J131").Value = 917.2917
$ws.Range("K131").Value = 2347.5
$ws.Range("L131").Value = 2751.8751
$ws.Range("M131").Value = 2692.5
$ws.Range("N131").Value = -12831.8751

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 12697.889
$ws.Range("J123").Value = 12697.889
$ws.Range("L123").Value = 12697.889
$ws.Range("N123").Value = -17597.889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 33985.645
$ws.Range("J127").Value = 33985.645
$ws.Range("L127").Value = 33985.645
$ws.Range("N127").Value = -43905.645
$ws.Range("H132").Value = 4415.96
$ws.Range("I132").Value = 5857.2856
$ws.Range("J132").Value = 2581.5454
$ws.Range("K132").Value = 17571.8568
$ws.Range("L132").Value = 7744.6362
$ws.Range("M132").Value = -15041.8568
$ws.Range("N132").Value = -12804.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1127.6154
$ws.Range("I126").Value = 1018.2
$ws.Range("J126").Value = 1492.3334
$ws.Range("K126").Value = 3054.6
$ws.Range("L126").Value = 4477.0002
$ws.Range("M126").Value = -584.6000000000004
$ws.Range("N126").Value = -9417.0002
$ws.Range("H132").Value = 4417.375
$ws.Range("I132").Value = 4765.8237
$ws.Range("J132").Value = 3571.1428
$ws.Range("K132").Value = 14297.4711
$ws.Range("L132").Value = 10713.4284
$ws.Range("M132").Value = -11767.4711
$ws.Range("N132").Value = -15773.4284
